$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 17.07792282104492
$ws.Range("C3").Value = 16.32094383239746
$ws.Range("C4").Value = 16.08920097351074
$ws.Range("C5").Value = 15.9461498260498
$ws.Range("C6").Value = 15.87390899658203
